# Adding 10 Manhattan plots lecture and recitation
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Schedule" (sheet1) — Week / Module / Topic
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Schedule")

# Row 9 (week 8): module+topic collapse into "Open session, capstone prep"
$ws1.Range("B9").Value = "Open session, capstone prep"
$ws1.Range("C9").Value = "Open session, capstone prep"

# Row 12 (week 11): Manhattan plots merges with "making lots of plots at once"
$ws1.Range("C12").Value = "Manhattan plots and making lots of plots at once"

# Row 14 (week 13): Interactive plots -> ggplot extension packages and complexheatmap
$ws1.Range("C14").Value = "ggplot extension packages and complexheatmap"

# New row 15 (week 14): No class, Thanksgiving / Relaxing and eating
$ws1.Range("B15").Value = "No class, Thanksgiving"
$ws1.Range("C15").Value = "Relaxing and eating"

# Shift the remaining two rows (old week 15/16 capstone rows) down to 16/17
$ws1.Range("A16").Value = 15
$ws1.Range("B16").Value = "4: Putting it together"
$ws1.Range("C16").Value = "Capstone assignment open session"

$ws1.Range("A17").Value = 16
$ws1.Range("B17").Value = "4: Putting it together"
$ws1.Range("C17").Value = "Capstone assignment open session"

# Column B got wider to fit the new, longer topic text
$ws1.Columns.Item(2).ColumnWidth = 24

# ---------------------------------------------------------------------------
# Sheet "Schedule_date" (sheet2) — Week / Date / Module / Topic
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Schedule_date")

$ws2.Range("C9").Value = "Open session, capstone prep"
$ws2.Range("D9").Value = "Open session, capstone prep"

$ws2.Range("D12").Value = "Manhattan plots and making lots of plots at once"

$ws2.Range("D14").Value = "ggplot extension packages and complexheatmap"

$ws2.Range("C15").Value = "No class, Thanksgiving"
$ws2.Range("D15").Value = "Relaxing and eating"

$ws2.Range("D16").Value = "Capstone assignment open session"

$ws2.Range("B17").Value = 44900
$ws2.Range("D17").Value = "Capstone assignment open session"

# Column widths tighten up to plain (non bestFit) custom widths
$ws2.Columns.Item(2).ColumnWidth = 13
$ws2.Columns.Item(3).ColumnWidth = 22.5

# ---------------------------------------------------------------------------
# Selections / active sheet — author left off on "Schedule" (tab 1) this time
# ---------------------------------------------------------------------------
$ws2.Range("A1:D17").Select()
$ws1.Activate()
$ws1.Range("F23").Select()
